$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append one new row (row 7) with the trip/shipment record captured at
# 2025-05-01T14:19:45 (local Arabic-formatted timestamp), mirroring the
# existing rows' layout:
#   A: notes (empty)            E: trip type
#   B: companion name            F: vehicle
#   C: quantity                  G: organisation
#   D: camp                      H: time

# Column A is blank for this record. Force the Text format first so the
# cell is still materialised (an empty cell assignment alone is dropped),
# matching the other rows which always carry a (possibly empty) value in
# column A.
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = ""

$ws.Range("B7").Value = "أحمد شريم"

# Column C holds quantities as text in this sheet (e.g. "2323"), not
# numbers, so force Text format before assigning, otherwise Excel would
# coerce the numeric-looking string into a real number.
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "2323"

$ws.Range("D7").Value = "ايتا"
$ws.Range("E7").Value = "الرحلة 2"
$ws.Range("F7").Value = "C2"
$ws.Range("G7").Value = "NRC"
$ws.Range("H7").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٥:١٩:٤٥ م"
